# Weekly price-sheet update: a new (most-recent) observation is inserted
# as row 94, pushing the existing historical rows (94-141) down by one
# row (95-142). The dimension grows from A1:R141 to A1:R142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 94, shifting rows 94:141 down to 95:142.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the latest data point.
$ws.Cells.Item(94, 1).Value = 4
$ws.Cells.Item(94, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(94, 3).Value = "Los Lagos"
$ws.Cells.Item(94, 4).Value = 44523
$ws.Cells.Item(94, 5).Value = 10
$ws.Cells.Item(94, 6).Value = 100112039
$ws.Cells.Item(94, 7).Value = "Ciboulette"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 240
$ws.Cells.Item(94, 11).Value = 2500
$ws.Cells.Item(94, 12).Value = 2500
$ws.Cells.Item(94, 13).Value = 2500
$ws.Cells.Item(94, 14).Value = "`$/docena de atados"
$ws.Cells.Item(94, 15).Value = "Región Metropolitana"
$ws.Cells.Item(94, 16).Value = 833
$ws.Cells.Item(94, 17).Value = 3
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# Note: Rows.Insert() already carries the "Fecha" date-number format (style
# index 2) down into the new D94 cell, matching the rest of column D, so no
# extra style assignment is necessary here.
